# Update the training schedule numbers in row 4 (E4, G4, H4) and move the
# active cell selection to E5, matching the refreshed task data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 5
$ws.Range("G4").Value = -3
$ws.Range("H4").Value = 13

# Leave the selection on E5, as last saved in the updated workbook.
$ws.Range("E5").Select()
